$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: the phone-number validation message in H2 is kept, but the adjacent
# "combined rich-text" note cell (M2) is cleared out (it used to hold a
# two-run rich-text note about Email + PhoneNumber errors). Re-apply the
# existing wrap/top-align formatting so the cell keeps its look after the
# value is removed.
$ws.Range("M2").ClearContents()
$ws.Range("M2").WrapText = $true
$ws.Range("M2").VerticalAlignment = -4160

# Rows 3-6: the status column (H) used to read "Ok" -- it now reports
# "User EXISTED" (the same text that used to live in the note column M).
# The per-row note column (M) is cleared since its text now lives in H.
$ws.Range("H3").Value2 = "User EXISTED"
$ws.Range("H3").WrapText = $true
$ws.Range("H3").VerticalAlignment = -4160
$ws.Range("M3").ClearContents()
$ws.Range("M3").WrapText = $true
$ws.Range("M3").VerticalAlignment = -4160

$ws.Range("H4").Value2 = "User EXISTED"
$ws.Range("H4").WrapText = $true
$ws.Range("H4").VerticalAlignment = -4160
$ws.Range("M4").ClearContents()
$ws.Range("M4").WrapText = $true
$ws.Range("M4").VerticalAlignment = -4160

$ws.Range("H5").Value2 = "User EXISTED"
$ws.Range("H5").WrapText = $true
$ws.Range("H5").VerticalAlignment = -4160
$ws.Range("M5").ClearContents()
$ws.Range("M5").WrapText = $true
$ws.Range("M5").VerticalAlignment = -4160

$ws.Range("H6").Value2 = "User EXISTED"
$ws.Range("H6").WrapText = $true
$ws.Range("H6").VerticalAlignment = -4160
$ws.Range("M6").ClearContents()
$ws.Range("M6").WrapText = $true
$ws.Range("M6").VerticalAlignment = -4160

# Move the active selection to H8 (below the data), matching the saved view.
$ws.Range("H8").Select() | Out-Null
